# Generate Report for Handback
# Swaps the "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" and
# "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" rows' identity across the
# Overview / zh-cn / de-de sheets, and updates the status, handoff/handback
# file names and timestamps to reflect a completed handback.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Preserve hyperlink target URLs (they keep pointing at the same commit
# URLs as before; only which row/display text they are attached to
# changes) before removing the stale hyperlink objects.
$ovB2Url = $ws.Hyperlinks.Item(1).Address
$ovB3Url = $ws.Hyperlinks.Item(2).Address
if (-not $ovB2Url) { $ovB2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07756e0ee203ef41db1134c32bbd483299de7bee/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" }
if (-not $ovB3Url) { $ovB3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91250d2e96694af68355197f85f180030b33170a/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" }
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("B2").Value = "e2e\1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("G2").Value = "2016-08-31 07:44:02"

$ws.Range("A3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("B3").Value = "e2e\d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-31 07:42:16"

$ws.Hyperlinks.Add($ws.Range("B2"), $ovB2Url, "", "", "e2e\1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $ovB3Url, "", "", "e2e\d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md")

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$zhA2Url = $ws.Hyperlinks.Item(1).Address
$zhI2Url = $ws.Hyperlinks.Item(2).Address
$zhA3Url = $ws.Hyperlinks.Item(3).Address
$zhI3Url = $ws.Hyperlinks.Item(4).Address
if (-not $zhA2Url) { $zhA2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07756e0ee203ef41db1134c32bbd483299de7bee/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" }
if (-not $zhI2Url) { $zhI2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7e4da89eb5f61aab0ed89ac8bee8b462e06be7e1/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" }
if (-not $zhA3Url) { $zhA3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91250d2e96694af68355197f85f180030b33170a/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" }
if (-not $zhI3Url) { $zhI3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7e4da89eb5f61aab0ed89ac8bee8b462e06be7e1/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" }
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("G2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-31 07:43:50"
$ws.Range("I2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("J2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-31 07:44:34"

$ws.Range("A3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-31 07:41:58"
$ws.Range("I3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("J3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.zh-cn.xlf"
$ws.Range("P3").Value = ""

$ws.Hyperlinks.Add($ws.Range("A2"), $zhA2Url, "", "", "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $zhI2Url, "", "", "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $zhA3Url, "", "", "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $zhI3Url, "", "", "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$deA2Url = $ws.Hyperlinks.Item(1).Address
$deI2Url = $ws.Hyperlinks.Item(2).Address
$deA3Url = $ws.Hyperlinks.Item(3).Address
$deI3Url = $ws.Hyperlinks.Item(4).Address
if (-not $deA2Url) { $deA2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07756e0ee203ef41db1134c32bbd483299de7bee/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" }
if (-not $deI2Url) { $deI2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c36129d8b4b179d5cec0bf65abbc209e95188763/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md" }
if (-not $deA3Url) { $deA3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91250d2e96694af68355197f85f180030b33170a/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" }
if (-not $deI3Url) { $deI3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c36129d8b4b179d5cec0bf65abbc209e95188763/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md" }
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("G2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.de-de.xlf"
$ws.Range("H2").Value = "2016-08-31 07:44:02"
$ws.Range("I2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$ws.Range("J2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.de-de.xlf"
$ws.Range("K2").Value = "2016-08-31 07:44:52"

$ws.Range("A3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.de-de.xlf"
$ws.Range("H3").Value = "2016-08-31 07:42:16"
$ws.Range("I3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"
$ws.Range("J3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.de-de.xlf"
$ws.Range("P3").Value = ""

$ws.Hyperlinks.Add($ws.Range("A2"), $deA2Url, "", "", "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $deI2Url, "", "", "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $deA3Url, "", "", "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $deI3Url, "", "", "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md")
